$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '30.378.71'
Set-TextValue 'E2' '  +11.84%  '
Set-TextValue 'D3' '1.878.17'
Set-TextValue 'E3' '  +8.11%  '
Set-TextValue 'D4' '0.9931'
Set-TextValue 'E4' '  -0.44%  '
Set-TextValue 'D5' '250.30'
Set-TextValue 'E5' '  +4.09%  '
Set-TextValue 'D6' '0.9917'
Set-TextValue 'E6' '  -0.62%  '
Set-TextValue 'D7' '0.4969'
Set-TextValue 'E7' '  +3.63%  '
Set-TextValue 'D8' '44.95'
Set-TextValue 'E8' '  +8.81%  '
Set-TextValue 'D9' '0.2853'
Set-TextValue 'E9' '  +10.02%  '
Set-TextValue 'D10' '0.06546'
Set-TextValue 'D11' '1.861.78'
Set-TextValue 'E11' '  +7.97%  '
Set-TextValue 'D12' '17.02'
Set-TextValue 'E12' '  +5.75%  '
Set-TextValue 'D13' '0.07166'
Set-TextValue 'E13' '  +3.39%  '
Set-TextValue 'D14' '0.6682'
Set-TextValue 'E14' '  +10.96%  '
Set-TextValue 'D15' '86.26'
Set-TextValue 'E15' '  +12.27%  '
Set-TextValue 'D16' '4.808'
Set-TextValue 'E16' '  +8.37%  '
Set-TextValue 'D17' '30.332.71'
Set-TextValue 'E17' '  +11.86%  '
Set-TextValue 'E18' '  -0.30%  '
Set-TextValue 'D19' '0.000007497'
Set-TextValue 'E19' '  +6.17%  '
Set-TextValue 'D20' '12.58'
Set-TextValue 'E20' '  +10.32%  '
Set-TextValue 'D21' '0.9929'
Set-TextValue 'E21' '  -0.49%  '
Set-TextValue 'D22' '2.089.58'
Set-TextValue 'D23' '4.696'
Set-TextValue 'E23' '  +6.37%  '
Set-TextValue 'D24' '5.510'
Set-TextValue 'E24' '  +7.84%  '
Set-TextValue 'D25' '8.990'
Set-TextValue 'E25' '  +7.22%  '
Set-TextValue 'D26' '143.88'
Set-TextValue 'E26' '  +1.51%  '
Set-TextValue 'D27' '134.75'
Set-TextValue 'E27' '  +26.23%  '
Set-TextValue 'D28' '16.78'
Set-TextValue 'E28' '  +10.12%  '
Set-TextValue 'D29' '1.944'
Set-TextValue 'E29' '  +6.79%  '
Set-TextValue 'D30' '1.399'
Set-TextValue 'E30' '  +1.86%  '
Set-TextValue 'D31' '4.255'
Set-TextValue 'E31' '  +7.90%  '
Set-TextValue 'D32' '0.08619'
Set-TextValue 'E32' '  +8.67%  '
Set-TextValue 'D33' '3.899'
Set-TextValue 'E33' '  +6.29%  '
Set-TextValue 'D34' '0.05053'
Set-TextValue 'E34' '  +6.58%  '
Set-TextValue 'D35' '1.137'
Set-TextValue 'E35' '  +12.38%  '
Set-TextValue 'D36' '0.6858'
Set-TextValue 'E36' '  +11.23%  '
Set-TextValue 'D37' '2.682'
Set-TextValue 'E37' '  +3.44%  '
Set-TextValue 'D38' '2.318'
Set-TextValue 'E38' '  +14.89%  '
Set-TextValue 'D39' '2.757'
Set-TextValue 'E39' '  +8.89%  '
Set-TextValue 'D40' '0.9500'
Set-TextValue 'E40' '  +3.19%  '
Set-TextValue 'D41' '0.01626'
Set-TextValue 'E41' '  +9.19%  '
Set-TextValue 'D42' '6.137'
Set-TextValue 'E42' '  +7.95%  '
Set-TextValue 'D43' '103.28'
Set-TextValue 'E43' '  +4.60%  '
Set-TextValue 'D44' '0.9922'
Set-TextValue 'E44' '  -0.55%  '
Set-TextValue 'D45' '0.4161'
Set-TextValue 'E45' '  +8.85%  '
Set-TextValue 'D46' '7.452'
Set-TextValue 'E46' '  +8.86%  '
Set-TextValue 'D47' '0.1245'
Set-TextValue 'E47' '  +8.28%  '
Set-TextValue 'D48' '0.05631'
Set-TextValue 'E48' '  +5.33%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '8.328'
Set-TextValue 'E49' '  +6.36%  '
Set-TextValue 'B50' 'Elrond'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D50' '32.23'
Set-TextValue 'E50' '  +7.80%  '
Set-TextValue 'D51' '1.342'
Set-TextValue 'E51' '  +8.05%  '
